$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 4")
$ws.Activate()

# Insert a new row above row 17, shifting existing rows (17-31) down to (18-32)
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new reaction data
$ws.Range("A17").Value = "R16"
$ws.Range("B17").Value = "Al2O_g"
$ws.Range("C17").Value = -9.36
$ws.Range("D17").Value = 41956

# Update the active selection to match the saved workbook state
$ws.Range("F14").Select()
